# feat: add upload file and supabase
#
# - Shorten the first patient's full name.
# - Replace the sample "LinkPdf" URL (old wordpress manual PDF link) with a
#   Google Drive share link on all three data rows.
# - Move the active-cell selection to D8 (from C12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: NombreCompleto (B2) - shorten the name.
$ws.Range("B2").Value2 = "Alberto Fulanito"

# LinkPdf column (L) for all three rows - new Google Drive link.
$newLink = "https://drive.google.com/file/d/1sp1StYshn1Fio4LQF9ORlnkTC0AQKnBI/view"
$ws.Range("L2").Value2 = $newLink
$ws.Range("L3").Value2 = $newLink
$ws.Range("L4").Value2 = $newLink

# Update the active selection shown when the workbook is opened.
$ws.Range("D8").Select()
